$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (UNH)
$ws.Range("D2").Value = 333.94
$ws.Range("E2").Value = 50.9
$ws.Range("F2").Value = 1.28
$ws.Range("N2").Value = 52.47848103381103

# Row 3 (MET)
$ws.Range("D3").Value = 77.90000000000001
$ws.Range("E3").Value = 42.3
$ws.Range("F3").Value = 2.06
$ws.Range("N3").Value = 52.47848103381103

# Row 4 (AIG) - only MACRO_SCORE changes
$ws.Range("N4").Value = 52.47848103381103

# Row 5 (PRU)
$ws.Range("D5").Value = 110.51
$ws.Range("E5").Value = 64.3
$ws.Range("F5").Value = 2.35
$ws.Range("N5").Value = 52.47848103381103
